$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("A2")
Write-Host "color before: $($r.Font.Color)"
Write-Host "colorindex before: $($r.Font.ColorIndex)"
Write-Host "themecolor before: $($r.Font.ThemeColor)"
$r.Font.Size = 8
Write-Host "color after: $($r.Font.Color)"
